$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: shift column headers left (C<-D<-E<-C)
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Update data rows 2-21: C becomes the species text, E becomes 0.5
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 3).Value = "s__Pediococcus acidilactici"
    $ws.Cells.Item($r, 5).Value = 0.5
}
